$d = $word.ActiveDocument

# Replace "the manager making decision" -> "the librarian making decision"
# This also merges the three adjacent runs (split around "the highest")
# into a single run, matching the target document.
$d.Content.Find.Execute("the manager making decision", $false, $false, $false, $false, $false, $true, 1, $false, "the librarian making decision", 2) | Out-Null

# Replace "Manager can trust" -> "Librarian can trust"
$d.Content.Find.Execute("Manager can trust", $false, $false, $false, $false, $false, $true, 1, $false, "Librarian can trust", 2) | Out-Null
